$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output ("Before count: " + $ws.Cells.FormatConditions.Count)
$r1 = $ws.Range("A54:XFD81")
$fc1 = $r1.FormatConditions.Add(9, 7, '=NOT(ISERROR(SEARCH("PWM",A1)))')
Write-Output ("After add count: " + $ws.Cells.FormatConditions.Count)
